$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "87.307.41"
Set-TextValue $ws "E2" "  +3.35%  "
Set-TextValue $ws "D3" "3.251.02"
Set-TextValue $ws "E3" "  -1.24%  "
Set-TextValue $ws "E4" "  -0.08%  "
Set-TextValue $ws "D5" "210.86"
Set-TextValue $ws "E5" "  -3.34%  "
Set-TextValue $ws "D6" "625.28"
Set-TextValue $ws "E6" "  -1.21%  "
Set-TextValue $ws "D7" "0.367"
Set-TextValue $ws "E7" "  +14.09%  "
Set-TextValue $ws "D8" "0.688"
Set-TextValue $ws "E8" "  +16.27%  "
Set-TextValue $ws "D9" "1.00"
Set-TextValue $ws "E9" "  +0.02%  "
Set-TextValue $ws "D10" "3.245.20"
Set-TextValue $ws "E10" "  -1.60%  "
Set-TextValue $ws "E11" "  -3.65%  "
Set-TextValue $ws "D12" "0.180"
Set-TextValue $ws "E12" "  +8.10%  "
Set-TextValue $ws "E13" "  -8.16%  "
Set-TextValue $ws "D14" "3.862.15"
Set-TextValue $ws "E14" "  -0.86%  "
Set-TextValue $ws "D15" "33.73"
Set-TextValue $ws "E15" "  -1.42%  "
Set-TextValue $ws "D16" "5.32"
Set-TextValue $ws "E16" "  -2.13%  "
Set-TextValue $ws "D17" "87.538.24"
Set-TextValue $ws "E17" "  +3.79%  "
Set-TextValue $ws "D18" "3.257.19"
Set-TextValue $ws "E18" "  -1.15%  "
Set-TextValue $ws "E19" "  -2.59%  "
Set-TextValue $ws "D20" "13.95"
Set-TextValue $ws "E20" "  -4.32%  "
Set-TextValue $ws "D21" "432.31"
Set-TextValue $ws "E21" "  -0.12%  "
Set-TextValue $ws "E22" "  -3.94%  "
Set-TextValue $ws "E23" "  +1.56%  "
Set-TextValue $ws "D24" "7.27"
Set-TextValue $ws "E24" "  -2.77%  "
Set-TextValue $ws "D25" "12.43"
Set-TextValue $ws "E25" "  +2.30%  "
Set-TextValue $ws "E26" "  -7.20%  "
Set-TextValue $ws "D27" "3.386.72"
Set-TextValue $ws "E27" "  -1.80%  "
Set-TextValue $ws "D28" "76.06"
Set-TextValue $ws "E28" "  -2.67%  "
Set-TextValue $ws "E29" "  -3.01%  "
Set-TextValue $ws "E30" "  -0.08%  "
Set-TextValue $ws "D31" "0.175"
Set-TextValue $ws "E31" "  +8.45%  "
Set-TextValue $ws "E32" "  -0.03%  "
Set-TextValue $ws "D33" "8.68"
Set-TextValue $ws "E33" "  -6.15%  "
Set-TextValue $ws "D34" "546.02"
Set-TextValue $ws "E34" "  -8.93%  "
Set-TextValue $ws "D35" "1.38"
Set-TextValue $ws "E35" "  -12.18%  "
Set-TextValue $ws "D36" "1.94"
Set-TextValue $ws "E36" "  -4.45%  "
Set-TextValue $ws "D37" "6.93"
Set-TextValue $ws "E37" "  +8.74%  "
Set-TextValue $ws "D38" "0.135"
Set-TextValue $ws "E38" "  -11.76%  "
Set-TextValue $ws "E39" "  -3.78%  "
Set-TextValue $ws "E40" "  -0.22%  "
Set-TextValue $ws "D41" "21.66"
Set-TextValue $ws "E41" "  +3.38%  "
Set-TextValue $ws "D42" "0.390"
Set-TextValue $ws "E42" "  -5.59%  "
Set-TextValue $ws "E43" "  -3.19%  "
Set-TextValue $ws "B44" "USDe"
Set-TextValue $ws "C44" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D44" "1.00"
Set-TextValue $ws "E44" "  -0.04%  "
Set-TextValue $ws "B45" "dogwifhat"
Set-TextValue $ws "C45" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D45" "2.91"
Set-TextValue $ws "E45" "  -6.02%  "
Set-TextValue $ws "D46" "154.90"
Set-TextValue $ws "E46" "  -2.59%  "
Set-TextValue $ws "D47" "178.61"
Set-TextValue $ws "E47" "  -6.10%  "
Set-TextValue $ws "D48" "44.83"
Set-TextValue $ws "E48" "  -0.46%  "
Set-TextValue $ws "E49" "  -4.99%  "
Set-TextValue $ws "B50" "Filecoin"
Set-TextValue $ws "C50" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D50" "4.19"
Set-TextValue $ws "E50" "  -1.32%  "
Set-TextValue $ws "B51" "Stellar"
Set-TextValue $ws "C51" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D51" "0.124"
Set-TextValue $ws "E51" "  +10.13%  "
